$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns retain their original text-detection-friendly
# behaviour (values like "1.00" or "6.70" must stay literal text, matching the
# original inline-string cells, rather than being auto-converted to numbers).
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = '69.065.56'
$ws.Range("E2").Value = '  -2.94%  '
$ws.Range("D3").Value = '3.515.30'
$ws.Range("E3").Value = '  -4.78%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '577.90'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = '171.21'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("D7").Value = '3.508.59'
$ws.Range("E7").Value = '  -4.73%  '
$ws.Range("D8").Value = '0.607'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("E10").Value = '  -5.61%  '
$ws.Range("D11").Value = '6.70'
$ws.Range("E11").Value = '  -2.09%  '
$ws.Range("D12").Value = '0.581'
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("D13").Value = '46.87'
$ws.Range("E13").Value = '  -4.16%  '
$ws.Range("E14").Value = '  -4.33%  '
$ws.Range("D15").Value = '4.080.31'
$ws.Range("E15").Value = '  -4.81%  '
$ws.Range("D16").Value = '8.51'
$ws.Range("E16").Value = '  -5.31%  '
$ws.Range("D17").Value = '620.68'
$ws.Range("E17").Value = '  -8.32%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.072.18'
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.496.77'
$ws.Range("E19").Value = '  -5.05%  '
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '17.41'
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").Value = '11.15'
$ws.Range("E22").Value = '  -3.69%  '
$ws.Range("D23").Value = '0.883'
$ws.Range("E23").Value = '  -6.05%  '
$ws.Range("D24").Value = '15.91'
$ws.Range("E24").Value = '  -8.37%  '
$ws.Range("D25").Value = '97.30'
$ws.Range("E25").Value = '  -4.49%  '
$ws.Range("E26").Value = '  -4.27%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '2.63'
$ws.Range("E28").Value = '  -6.53%  '
$ws.Range("D29").Value = '9.33'
$ws.Range("E29").Value = '  -9.00%  '
$ws.Range("D30").Value = '32.55'
$ws.Range("E30").Value = '  -7.19%  '
$ws.Range("D31").Value = '3.16'
$ws.Range("E31").Value = '  -7.59%  '
$ws.Range("D32").Value = '8.54'
$ws.Range("E32").Value = '  -6.89%  '
$ws.Range("E33").Value = '  -7.27%  '
$ws.Range("D34").Value = '6.97'
$ws.Range("E34").Value = '  -6.77%  '
$ws.Range("D35").Value = '634.03'
$ws.Range("E35").Value = '  +8.02%  '
$ws.Range("D36").Value = '10.74'
$ws.Range("E36").Value = '  -3.90%  '
$ws.Range("E37").Value = '  -5.27%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '56.70'
$ws.Range("E38").Value = '  -3.82%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = '3.40'
$ws.Range("E39").Value = '  -16.62%  '
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("D41").Value = '0.0446'
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("E42").Value = '  -5.74%  '
$ws.Range("D43").Value = '3.380.87'
$ws.Range("E43").Value = '  -8.13%  '
$ws.Range("E44").Value = '  -6.34%  '
$ws.Range("D45").Value = '32.85'
$ws.Range("E45").Value = '  -7.24%  '
$ws.Range("D46").Value = '0.0₃0688'
$ws.Range("E46").Value = '  -9.88%  '
$ws.Range("D47").Value = '2.56'
$ws.Range("E47").Value = '  -7.13%  '
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -4.44%  '
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '5.72'
$ws.Range("E50").Value = '  +15.02%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '132.22'
$ws.Range("E51").Value = '  -2.47%  '

# Restore the default (unstyled) cell style so no stray formatting is introduced.
$numRng.Style = "Normal"
